$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "Cantidad de 0s =" label with COUNTIF(...,1) formulas across I:Q
$ws.Range("A9").Value = "Cantidad de 0s ="
$ws.Range("I9").Formula = "=COUNTIF(I1:I7,1)"
$ws.Range("J9").Formula = "=COUNTIF(J1:J7,1)"
$ws.Range("K9").Formula = "=COUNTIF(K1:K7,1)"
$ws.Range("L9").Formula = "=COUNTIF(L1:L7,1)"
$ws.Range("M9").Formula = "=COUNTIF(M1:M7,1)"
$ws.Range("N9").Formula = "=COUNTIF(N1:N7,1)"
$ws.Range("O9").Formula = "=COUNTIF(O1:O7,1)"
$ws.Range("P9").Formula = "=COUNTIF(P1:P7,1)"
$ws.Range("Q9").Formula = "=COUNTIF(Q1:Q7,1)"

# Row 10: "Cantidad de 1s =" label with COUNTIF(...,0) formulas across I:Q
$ws.Range("A10").Value = "Cantidad de 1s ="
$ws.Range("I10").Formula = "=COUNTIF(I1:I7,0)"
$ws.Range("J10").Formula = "=COUNTIF(J1:J7,0)"
$ws.Range("K10").Formula = "=COUNTIF(K1:K7,0)"
$ws.Range("L10").Formula = "=COUNTIF(L1:L7,0)"
$ws.Range("M10").Formula = "=COUNTIF(M1:M7,0)"
$ws.Range("N10").Formula = "=COUNTIF(N1:N7,0)"
$ws.Range("O10").Formula = "=COUNTIF(O1:O7,0)"
$ws.Range("P10").Formula = "=COUNTIF(P1:P7,0)"
$ws.Range("Q10").Formula = "=COUNTIF(Q1:Q7,0)"
